$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes literal text into a cell without Excel auto-converting
# percent-looking strings ("51%") into a numeric percentage value. We stage
# the text in a scratch cell (forced to Text format), copy it, and paste
# *values only* into the destination so the destination keeps its original
# style/number format untouched.
function Set-LiteralText($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range("E2").Value = "2026-02-26 23:18:27"
Set-LiteralText "H2" "51%"
$ws.Range("K2").Value = "13.2 MJ/m2"
$ws.Range("O2").Value = "5.4 °C"
$ws.Range("E3").Value = "2026-02-26 23:18:29"
$ws.Range("E4").Value = "2026-02-26 23:18:31"
$ws.Range("J4").Value = "1027.2 hPa"
$ws.Range("O4").Value = "10.4 °C"
$ws.Range("E5").Value = "2026-02-26 23:18:34"
$ws.Range("O5").Value = "5.0 °C"
$ws.Range("E6").Value = "2026-02-26 23:18:37"
$ws.Range("E7").Value = "2026-02-26 23:18:39"
Set-LiteralText "H7" "83%"
$ws.Range("O7").Value = "12.9 °C"
$ws.Range("E8").Value = "2026-02-26 23:18:42"
$ws.Range("E9").Value = "2026-02-26 23:18:45"
$ws.Range("E10").Value = "2026-02-26 23:18:46"
$ws.Range("O10").Value = "9.1 °C"
$ws.Range("E11").Value = "2026-02-26 23:18:47"
$ws.Range("E12").Value = "2026-02-26 23:18:48"
$ws.Range("N12").Value = "6.6 °C 22:55 TU"
$ws.Range("O12").Value = "11.0 °C"
$ws.Range("E13").Value = "2026-02-26 23:18:49"
$ws.Range("J13").Value = "1028.4 hPa"
$ws.Range("O13").Value = "6.9 °C"
$ws.Range("E14").Value = "2026-02-26 23:18:50"
$ws.Range("O14").Value = "11.1 °C"
$ws.Range("E15").Value = "2026-02-26 23:18:52"
$ws.Range("O15").Value = "11.2 °C"
$ws.Range("E16").Value = "2026-02-26 23:18:53"
$ws.Range("E17").Value = "2026-02-26 23:18:54"
$ws.Range("K17").Value = "17.8 MJ/m2"
$ws.Range("E18").Value = "2026-02-26 23:18:55"
$ws.Range("O18").Value = "11.8 °C"
$ws.Range("E19").Value = "2026-02-26 23:18:56"
Set-LiteralText "H19" "50%"
$ws.Range("O19").Value = "11.2 °C"
$ws.Range("E20").Value = "2026-02-26 23:18:59"
Set-LiteralText "H20" "49%"
$ws.Range("O20").Value = "2.4 °C"
$ws.Range("E21").Value = "2026-02-26 23:19:01"
$ws.Range("E22").Value = "2026-02-26 23:19:04"
$ws.Range("E23").Value = "2026-02-26 23:19:07"
Set-LiteralText "H23" "40%"
$ws.Range("E24").Value = "2026-02-26 23:19:09"
Set-LiteralText "H24" "76%"
$ws.Range("O24").Value = "10.1 °C"
$ws.Range("E25").Value = "2026-02-26 23:19:12"
$ws.Range("E26").Value = "2026-02-26 23:19:15"
$ws.Range("J26").Value = "1024.6 hPa"
$ws.Range("O26").Value = "10.6 °C"
$ws.Range("E27").Value = "2026-02-26 23:19:18"
Set-LiteralText "H27" "42%"
$ws.Range("K27").Value = "17.2 MJ/m2"
$ws.Range("E28").Value = "2026-02-26 23:19:20"
$ws.Range("N28").Value = "4.9 °C 22:30 TU"
$ws.Range("O28").Value = "10.4 °C"
$ws.Range("E29").Value = "2026-02-26 23:19:23"
$ws.Range("N29").Value = "6.6 °C 22:59 TU"
$ws.Range("O29").Value = "11.2 °C"
$ws.Range("E30").Value = "2026-02-26 23:19:26"
$ws.Range("E31").Value = "2026-02-26 23:19:28"
$ws.Range("E32").Value = "2026-02-26 23:19:31"
Set-LiteralText "H32" "69%"
$ws.Range("O32").Value = "7.4 °C"
$ws.Range("E33").Value = "2026-02-26 23:19:34"
$ws.Range("O33").Value = "8.4 °C"
$ws.Range("E34").Value = "2026-02-26 23:19:36"
$ws.Range("O34").Value = "4.5 °C"
$ws.Range("E35").Value = "2026-02-26 23:19:39"
$ws.Range("E36").Value = "2026-02-26 23:19:42"
$ws.Range("E37").Value = "2026-02-26 23:19:44"
Set-LiteralText "H37" "76%"
$ws.Range("O37").Value = "7.4 °C"
$ws.Range("E38").Value = "2026-02-26 23:19:47"
$ws.Range("E39").Value = "2026-02-26 23:19:49"
Set-LiteralText "H39" "41%"
$ws.Range("E40").Value = "2026-02-26 23:19:52"
$ws.Range("J40").Value = "1027.6 hPa"
$ws.Range("O40").Value = "9.2 °C"
$ws.Range("E41").Value = "2026-02-26 23:19:55"
$ws.Range("E42").Value = "2026-02-26 23:19:57"
Set-LiteralText "H42" "89%"
$ws.Range("N42").Value = "6.4 °C 22:56 TU"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-26 23:20:00"
$ws.Range("E44").Value = "2026-02-26 23:20:02"
$ws.Range("K44").Value = "16.6 MJ/m2"
$ws.Range("O44").Value = "2.1 °C"
$ws.Range("E45").Value = "2026-02-26 23:20:05"
$ws.Range("O45").Value = "10.3 °C"
$ws.Range("E46").Value = "2026-02-26 23:20:08"
